$d = $word.ActiveDocument

# Apply an 11pt (sz=22 half-points) font size to every paragraph (including
# paragraph-mark run properties) and every run's character formatting,
# matching both the regular (sz) and complex-script (szCs) sizes.
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Size = 11
    $p.Range.Font.SizeBi = 11
}

# Move the "_GoBack" bookmark from the end of the document to the start of
# the "Any other supporting documents..." list item.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Any other supporting documents as specified in the instructions*") {
        $target = $p
    }
}

$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$r = $d.Range($target.Range.Start, $target.Range.Start)
$d.Bookmarks.Add("_GoBack", $r)
